$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 62157.41
$ws.Range("I64").Value = 145570.86
$ws.Range("J64").Value = 3768
$ws.Range("K64").Value = 145570.86
$ws.Range("L64").Value = 3768
$ws.Range("M64").Value = -145322.86
$ws.Range("N64").Value = -4264

$ws.Range("H67").Value = 62157.41
$ws.Range("I67").Value = 145570.86
$ws.Range("J67").Value = 3768
$ws.Range("K67").Value = 145570.86
$ws.Range("L67").Value = 3768
$ws.Range("M67").Value = -144712.86
$ws.Range("N67").Value = -5484

$ws.Range("H80").Value = 53216.58
$ws.Range("I80").Value = 313.8889
$ws.Range("J80").Value = 100829
$ws.Range("K80").Value = 941.6667
$ws.Range("L80").Value = 302487
$ws.Range("M80").Value = 56.33330000000001
$ws.Range("N80").Value = -304483

$ws.Range("H83").Value = 53216.58
$ws.Range("I83").Value = 313.8889
$ws.Range("J83").Value = 100829
$ws.Range("K83").Value = 2825.0001
$ws.Range("L83").Value = 907461
$ws.Range("M83").Value = 2166.9999
$ws.Range("N83").Value = -917445

$ws.Range("H88").Value = 4550.25
$ws.Range("I88").Value = 2233.3333
$ws.Range("J88").Value = 5322.5557
$ws.Range("K88").Value = 2233.3333
$ws.Range("L88").Value = 5322.5557
$ws.Range("M88").Value = -1827.3333
$ws.Range("N88").Value = -6134.5557

$ws.Range("H91").Value = 4550.25
$ws.Range("I91").Value = 2233.3333
$ws.Range("J91").Value = 5322.5557
$ws.Range("K91").Value = 2233.3333
$ws.Range("L91").Value = 5322.5557
$ws.Range("M91").Value = -829.3332999999998
$ws.Range("N91").Value = -8130.5557

$ws.Range("H98").Value = 762.5714
$ws.Range("I98").Value = 733.0333000000001
$ws.Range("J98").Value = 939.8
$ws.Range("K98").Value = 733.0333000000001
$ws.Range("L98").Value = 939.8
$ws.Range("M98").Value = 764.9666999999999
$ws.Range("N98").Value = -3935.8

$ws.Range("H122").Value = 762.5714
$ws.Range("I122").Value = 733.0333000000001
$ws.Range("J122").Value = 939.8
$ws.Range("K122").Value = 2199.0999
$ws.Range("L122").Value = 2819.4
$ws.Range("M122").Value = 250.9000999999998
$ws.Range("N122").Value = -7719.4

$ws.Range("H132").Value = 879.3333
$ws.Range("I132").Value = 873.1111
$ws.Range("K132").Value = 2619.3333
$ws.Range("M132").Value = -89.33329999999978

$ws.Range("H137").Value = 1370.82
$ws.Range("I137").Value = 896.0714
$ws.Range("J137").Value = 1555.4445
$ws.Range("K137").Value = 2688.2142
$ws.Range("L137").Value = 4666.333500000001
$ws.Range("M137").Value = -138.2142000000003
$ws.Range("N137").Value = -9766.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 68205.13
$ws.Range("I2").Value = 1423.0834
$ws.Range("K2").Value = 1423.0834
$ws.Range("M2").Value = -1310.0834

$ws.Range("H32").Value = 26054.855
$ws.Range("I32").Value = 4752.1875
$ws.Range("K32").Value = 4752.1875
$ws.Range("M32").Value = -4465.1875

$ws.Range("H63").Value = 4000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 4000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -26864

$ws.Range("H116").Value = 68205.13
$ws.Range("I116").Value = 1423.0834
$ws.Range("K116").Value = 1423.0834
$ws.Range("M116").Value = 870.9166

$ws.Range("H122").Value = 2331.238
$ws.Range("I122").Value = 2251.3845
$ws.Range("J122").Value = 2461
$ws.Range("K122").Value = 6754.1535
$ws.Range("L122").Value = 7383
$ws.Range("M122").Value = -4304.1535
$ws.Range("N122").Value = -12283

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 68205.13
$ws.Range("I3").Value = 1423.0834
$ws.Range("K3").Value = 1423.0834
$ws.Range("M3").Value = -1309.0834

$ws.Range("H35").Value = 19083.818
$ws.Range("J35").Value = 19992.2
$ws.Range("L35").Value = 19992.2
$ws.Range("N35").Value = -20612.2

$ws.Range("H94").Value = 27361.947
$ws.Range("I94").Value = 125655.5
$ws.Range("J94").Value = 1150.3334
$ws.Range("K94").Value = 125655.5
$ws.Range("L94").Value = 1150.3334
$ws.Range("M94").Value = -125204.5
$ws.Range("N94").Value = -2052.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 438.15384
$ws.Range("I22").Value = 299
$ws.Range("J22").Value = 557.4286
$ws.Range("K22").Value = 299
$ws.Range("L22").Value = 557.4286
$ws.Range("M22").Value = 51
$ws.Range("N22").Value = -1257.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 16766.6
$ws.Range("I87").Value = 4166.3335
$ws.Range("J87").Value = 19916.666
$ws.Range("K87").Value = 12499.0005
$ws.Range("L87").Value = 59749.99800000001
$ws.Range("M87").Value = -11251.0005
$ws.Range("N87").Value = -62245.99800000001

$ws.Range("H90").Value = 16766.6
$ws.Range("I90").Value = 4166.3335
$ws.Range("J90").Value = 19916.666
$ws.Range("K90").Value = 37497.0015
$ws.Range("L90").Value = 179249.994
$ws.Range("M90").Value = -31257.0015
$ws.Range("N90").Value = -191729.994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 186544.1
$ws.Range("I70").Value = 254498.12
$ws.Range("J70").Value = 5333.3335
$ws.Range("K70").Value = 254498.12
$ws.Range("L70").Value = 5333.3335
$ws.Range("M70").Value = -254228.12
$ws.Range("N70").Value = -5873.3335

$ws.Range("H73").Value = 186544.1
$ws.Range("I73").Value = 254498.12
$ws.Range("J73").Value = 5333.3335
$ws.Range("K73").Value = 254498.12
$ws.Range("L73").Value = 5333.3335
$ws.Range("M73").Value = -253562.12
$ws.Range("N73").Value = -7205.3335

$ws.Range("H80").Value = 3050
$ws.Range("I80").Value = 2924
$ws.Range("J80").Value = 3260
$ws.Range("K80").Value = 2924
$ws.Range("L80").Value = 3260
$ws.Range("M80").Value = -1926
$ws.Range("N80").Value = -5256

$ws.Range("H83").Value = 3050
$ws.Range("I83").Value = 2924
$ws.Range("J83").Value = 3260
$ws.Range("K83").Value = 14620
$ws.Range("L83").Value = 16300
$ws.Range("M83").Value = -9628
$ws.Range("N83").Value = -26284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3266.6667
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 3612.5
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 3612.5
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -3988.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 987.625
$ws.Range("I122").Value = 984.2
$ws.Range("K122").Value = 2952.6
$ws.Range("M122").Value = -502.6000000000004
